# The sheet contains a daily price log for "Pepino ensalada" at
# "Feria Lagunitas de Puerto Montt". The commit adds a new weekly
# observation. In the OOXML diff this shows up as: a brand-new row is
# inserted right after row 242, pushing the former rows 243-340 down by
# one (to 244-341), and the newly inserted row 243 carries the new
# reading (date 2022-11-01 -> serial 44875, volume 200, min/max/avg price
# 25000, and $/Kg price 417). All the other (constant) columns simply
# come along with the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 243:340 down to 244:341, creating a blank row 243.
$ws.Rows.Item(243).Insert()

# Seed the new row 243 with the same constant columns as its neighbour
# (which now lives in row 244), then overwrite the columns that actually
# differ for the new observation.
$ws.Range("A243:R243").Value2 = $ws.Range("A244:R244").Value2

$ws.Range("D243").Value2 = 44875
$ws.Range("J243").Value2 = 200
$ws.Range("K243").Value2 = 25000
$ws.Range("L243").Value2 = 25000
$ws.Range("M243").Value2 = 25000
$ws.Range("P243").Value2 = 417
